# Global_Parameters.xlsx - add "Mip Gap" solver-option section
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert three new rows above the current "Scaling" section (old row 9).
#    This pushes the old rows 9-19 down to 12-22.
$ws.Rows("9:11").Insert()

# 2) Fix up cell styles for the new rows by copying formats from existing,
#    similarly-styled cells elsewhere on the sheet.
$ws.Range("B5").Copy()
$ws.Range("C9").PasteSpecial(-4122)

$ws.Range("B7:C7").Copy()
$ws.Range("B10:C10").PasteSpecial(-4122)

$ws.Range("C15").Copy()
$ws.Range("C11").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 3) Populate the new "Mip Gap" rows
$ws.Range("B10").Value = "Mip Gap"
$ws.Range("C10").Value = "[%]"

$ws.Range("B11").Value = "pMIPGap"
$ws.Range("C11").Value = 0.05
$ws.Range("E11").Value = "Relative MIP gap"
$ws.Range("F11").Value = "The MIP solver will terminate (with an optimal result) when the gap between the lower and upper objective bound is less than pMIPGap"
$ws.Range("G11").Value = "Factor"
$ws.Range("H11").Value = 0.05

# 4) Conditional formatting: the old rules for the shifted cells point at the
#    wrong rows now (still "C12" / "C19" instead of "C15" / "C22"), and the
#    new pMIPGap cell (C11) needs the same Yes/No colouring the other
#    boolean-ish value cells have.
$ws.Range("C12").FormatConditions.Delete()
$ws.Range("C19").FormatConditions.Delete()

function Add-YesNoFormat($addr) {
    $rng = $ws.Range($addr)
    $fcNo = $rng.FormatConditions.Add(1, 3, "=""No""")
    $fcNo.Font.Bold = $true
    $fcNo.Font.Italic = $false
    $fcNo.Font.Color = 4824142
    $fcYes = $rng.FormatConditions.Add(1, 3, "=""Yes""")
    $fcYes.Font.Bold = $true
    $fcYes.Font.Italic = $false
    $fcYes.Font.Color = 3473849
}

Add-YesNoFormat "C11"
Add-YesNoFormat "C15"
Add-YesNoFormat "C22"

# 5) Data validations: extend / move sqrefs to match the new layout.
$ws.Range("C8:C9").Validation.Delete()
$ws.Range("C8:C9").Validation.Add(3, 1, 1, "No, Yes")
$ws.Range("C8:C9").Validation.ShowInput = $true
$ws.Range("C8:C9").Validation.ShowError = $true

$ws.Range("C5").Validation.Delete()
$ws.Range("C15").Validation.Delete()
$ws.Range("C18").Validation.Delete()
$ws.Range("C11").Validation.Delete()
$ws.Range("C5,C15,C18,C11").Validation.Add(0)
$ws.Range("C5,C15,C18,C11").Validation.ShowInput = $true
$ws.Range("C5,C15,C18,C11").Validation.ShowError = $true

$ws.Range("C22").Validation.Delete()
$ws.Range("C22").Validation.Add(1, 4, 5, "0")
$ws.Range("C22").Validation.ShowInput = $true
$ws.Range("C22").Validation.ShowError = $true
